$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the bogus/placeholder row ("FUNDAÇÃO GETULIO VARGAS" / "123.456.789-10")
# that had been left in the data at row 98. Deleting the entire row shifts every
# following row up by one (old row 99 "CENTRO DE ENSINO MEDIO ARY RIBEIRO
# VALADAO FILHO" becomes the new row 98, ..., old row 152 becomes the new row
# 151) and the two now-unused strings drop out of the shared-string table.
$ws.Rows.Item(98).Delete()

# Leave the selection where the author ended up after the edit.
$ws.Range("C94").Select()
